# Applies two corrections to the transcription:
#   1. "la S<exp>ain</exp>t Jehan" -> "la S<exp>ainc</exp>t Jehan"
#      (plain text replace; the "<exp>" markers are literal characters here,
#       not real XML)
#   2. The lone "l" in " ou quave" + "l" + "q le " becomes "c" (forming
#      "quavecq"), and that run's explicit black color (w:color val="000000")
#      is dropped so only <w:rtl val="0"/> remains in its rPr - matching a
#      run elsewhere in the document that never had an explicit color.

$d = $word.ActiveDocument

# --- Change 1: "ain" -> "ainc" inside the "S<exp>...</exp>t" abbreviation ---
$d.Content.Find.Execute(
    "la S<exp>ain</exp>t Jehan",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "la S<exp>ainc</exp>t Jehan",
    2
) | Out-Null

# --- Change 2: single-letter run "l" -> "c", with its w:color removed ---

# Locate a pristine run elsewhere in the document that already carries no
# explicit <w:color> (only <w:rtl val="0"/>), so we can clone its (lack of)
# direct color formatting onto our target run without Word's "merge runs
# with identical resulting formatting" auto-coalescing kicking in.
$refAnchor = $d.Content
$refAnchor.Find.Execute(
    "couppe a la facon quo",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$colorlessRef = $d.Range($refAnchor.End, $refAnchor.End + 1)

# Locate the target single-character run (the "l" in " ou quave" / "l" / "q le ").
$target = $d.Content
$target.Find.Execute(
    "quave",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$targetRun = $d.Range($target.End, $target.End + 1)

# Copy the colorless run's formatted text (text + run formatting) onto the
# target run, which carries over the "no explicit color" property set.
$targetRun.FormattedText = $colorlessRef.FormattedText

# Now update the character itself from "n" (borrowed) to the required "c".
$targetRun2 = $d.Range($target.End, $target.End + 1)
$targetRun2.Text = "c"
